$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.225296442687747
$ws.Range("C2").Value = 0.5019762845849802
$ws.Range("P2").Value = 0.1739130434782609
$ws.Range("S2").Value = 0.09881422924901186
$ws.Range("B3").Value = 0.01587301587301587
$ws.Range("C3").Value = 0.01587301587301587
$ws.Range("J3").Value = 0.01587301587301587
$ws.Range("P3").Value = 0.753968253968254
$ws.Range("S3").Value = 0.1984126984126984
$ws.Range("J4").Value = 0.103448275862069
$ws.Range("P4").Value = 0.4827586206896552
$ws.Range("S4").Value = 0.4137931034482759
$ws.Range("B6").Value = 0.05725190839694656
$ws.Range("F6").Value = 0.07251908396946564
$ws.Range("J6").Value = 0.2442748091603053
$ws.Range("O6").Value = 0.01526717557251908
$ws.Range("Q6").Value = 0.1145038167938931
$ws.Range("R6").Value = 0.08015267175572519
$ws.Range("S6").Value = 0.4160305343511451
$ws.Range("B7").Value = 0.06511627906976744
$ws.Range("D7").Value = 0.02325581395348837
$ws.Range("E7").Value = 0.004651162790697674
$ws.Range("F7").Value = 0.06046511627906977
$ws.Range("J7").Value = 0.1209302325581395
$ws.Range("O7").Value = 0.02325581395348837
$ws.Range("Q7").Value = 0.1395348837209302
$ws.Range("R7").Value = 0.08372093023255814
$ws.Range("S7").Value = 0.4790697674418605
$ws.Range("B8").Value = 0.08333333333333333
$ws.Range("D8").Value = 0.01096491228070175
$ws.Range("F8").Value = 0.05263157894736842
$ws.Range("J8").Value = 0.1096491228070175
$ws.Range("O8").Value = 0.01754385964912281
$ws.Range("Q8").Value = 0.1513157894736842
$ws.Range("R8").Value = 0.08991228070175439
$ws.Range("S8").Value = 0.4846491228070176
$ws.Range("B9").Value = 0.08053691275167785
$ws.Range("D9").Value = 0.01677852348993289
$ws.Range("F9").Value = 0.05704697986577181
$ws.Range("J9").Value = 0.09731543624161074
$ws.Range("O9").Value = 0.01677852348993289
$ws.Range("Q9").Value = 0.1375838926174497
$ws.Range("R9").Value = 0.1174496644295302
$ws.Range("S9").Value = 0.4765100671140939
$ws.Range("B10").Value = 0.07412898443291327
$ws.Range("D10").Value = 0.01408450704225352
$ws.Range("F10").Value = 0.08080059303187546
$ws.Range("J10").Value = 0.09117865085248332
$ws.Range("O10").Value = 0.01630837657524092
$ws.Range("Q10").Value = 0.1830985915492958
$ws.Range("R10").Value = 0.09710896960711639
$ws.Range("S10").Value = 0.4432913269088213
$ws.Range("G11").Value = 0.1339869281045752
$ws.Range("J11").Value = 0.09477124183006536
$ws.Range("K11").Value = 0.2124183006535948
$ws.Range("L11").Value = 0.5359477124183006
$ws.Range("S11").Value = 0.02287581699346405
$ws.Range("G12").Value = 0.8
$ws.Range("J12").Value = 0.1257142857142857
$ws.Range("K12").Value = 0.005714285714285714
$ws.Range("L12").Value = 0.04
$ws.Range("S12").Value = 0.02857142857142857
$ws.Range("G13").Value = 0.7346938775510204
$ws.Range("J13").Value = 0.2448979591836735
$ws.Range("S13").Value = 0.02040816326530612
$ws.Range("F15").Value = 0.02293577981651376
$ws.Range("H15").Value = 0.1100917431192661
$ws.Range("I15").Value = 0.1284403669724771
$ws.Range("J15").Value = 0.3440366972477064
$ws.Range("K15").Value = 0.06422018348623854
$ws.Range("M15").Value = 0.009174311926605505
$ws.Range("N15").Value = 0.004587155963302753
$ws.Range("O15").Value = 0.04128440366972477
$ws.Range("S15").Value = 0.2752293577981652
$ws.Range("F16").Value = 0.01333333333333333
$ws.Range("H16").Value = 0.1533333333333333
$ws.Range("I16").Value = 0.1133333333333333
$ws.Range("J16").Value = 0.4466666666666667
$ws.Range("K16").Value = 0.08666666666666667
$ws.Range("M16").Value = 0.006666666666666667
$ws.Range("N16").Value = 0.006666666666666667
$ws.Range("O16").Value = 0.04
$ws.Range("S16").Value = 0.1333333333333333
$ws.Range("F17").Value = 0.02179176755447942
$ws.Range("H17").Value = 0.1694915254237288
$ws.Range("I17").Value = 0.12590799031477
$ws.Range("J17").Value = 0.423728813559322
$ws.Range("K17").Value = 0.06053268765133172
$ws.Range("M17").Value = 0.01452784503631961
$ws.Range("O17").Value = 0.04116222760290557
$ws.Range("S17").Value = 0.1428571428571428
$ws.Range("F18").Value = 0.02419354838709677
$ws.Range("H18").Value = 0.125
$ws.Range("I18").Value = 0.1411290322580645
$ws.Range("J18").Value = 0.4677419354838709
$ws.Range("K18").Value = 0.05241935483870968
$ws.Range("M18").Value = 0.01209677419354839
$ws.Range("O18").Value = 0.03225806451612903
$ws.Range("S18").Value = 0.1451612903225807
$ws.Range("F19").Value = 0.01948460087994972
$ws.Range("H19").Value = 0.1961030798240101
$ws.Range("I19").Value = 0.1043368950345695
$ws.Range("J19").Value = 0.3614079195474544
$ws.Range("K19").Value = 0.1081081081081081
$ws.Range("M19").Value = 0.02262727844123193
$ws.Range("N19").Value = 0.00251414204902577
$ws.Range("O19").Value = 0.04128440366972477
$ws.Range("S19").Value = 0.2752293577981652
